# Updated CHE_grids model - 2025-08-09 16:56
#
# The "solar" sheet's AG column ("grid_cell") values get reshuffled to a new
# CHE_<n> assignment per row (rows 4-26). Rows 19-22 keep their original
# values; all other rows in that range change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$gridCellByRow = @{
    4  = "CHE_7"
    5  = "CHE_12"
    6  = "CHE_14"
    7  = "CHE_18"
    8  = "CHE_11"
    9  = "CHE_15"
    10 = "CHE_25"
    11 = "CHE_3"
    12 = "CHE_13"
    13 = "CHE_24"
    14 = "CHE_5"
    15 = "CHE_8"
    16 = "CHE_21"
    17 = "CHE_9"
    18 = "CHE_4"
    19 = "CHE_20"
    20 = "CHE_1"
    21 = "CHE_6"
    22 = "CHE_0"
    23 = "CHE_17"
    24 = "CHE_19"
    25 = "CHE_10"
    26 = "CHE_22"
}

foreach ($row in $gridCellByRow.Keys) {
    $ws.Cells.Item($row, 33).Value = $gridCellByRow[$row]
}
